# Update Leve Profit sheets with refreshed Market Board price data
# (matches scheduled-runner market data refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value2 = 930.3333
$ws.Range("I18").Value2 = 930.3333
$ws.Range("K18").Value2 = 930.3333
$ws.Range("M18").Value2 = -646.3333

# Row 33
$ws.Range("H33").Value2 = 443.625
$ws.Range("I33").Value2 = 364.2857
$ws.Range("K33").Value2 = 364.2857
$ws.Range("M33").Value2 = -135.2857

# Row 123
$ws.Range("H123").Value2 = 26838.666
$ws.Range("J123").Value2 = 26838.666
$ws.Range("L123").Value2 = 26838.666
$ws.Range("N123").Value2 = -36638.666

# Row 132
$ws.Range("H132").Value2 = 3862.4324
$ws.Range("I132").Value2 = 4215.5312
$ws.Range("K132").Value2 = 12646.5936
$ws.Range("M132").Value2 = -10116.5936

# Row 134
$ws.Range("H134").Value2 = 59775
$ws.Range("J134").Value2 = 59775
$ws.Range("L134").Value2 = 59775
$ws.Range("N134").Value2 = -69915

# Row 135
$ws.Range("H135").Value2 = 14723721
$ws.Range("I135").Value2 = 462.72415
$ws.Range("J135").Value2 = 75720080
$ws.Range("K135").Value2 = 4164.51735
$ws.Range("L135").Value2 = 681480720
$ws.Range("M135").Value2 = -1629.51735
$ws.Range("N135").Value2 = -681485790

# Row 137
$ws.Range("H137").Value2 = 45456836
$ws.Range("I137").Value2 = 1666.4
$ws.Range("J137").Value2 = 83336140
$ws.Range("K137").Value2 = 4999.200000000001
$ws.Range("L137").Value2 = 250008420
$ws.Range("M137").Value2 = -2449.200000000001
$ws.Range("N137").Value2 = -250013520

# Row 138
$ws.Range("H138").Value2 = 2512.2957
$ws.Range("I138").Value2 = 2210.2
$ws.Range("J138").Value2 = 2733.3416
$ws.Range("K138").Value2 = 6630.599999999999
$ws.Range("L138").Value2 = 8200.024800000001
$ws.Range("M138").Value2 = -1490.599999999999
$ws.Range("N138").Value2 = -18480.0248

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 8157.2
$ws.Range("I32").Value2 = 5705.351
$ws.Range("J32").Value2 = 18907.615
$ws.Range("K32").Value2 = 5705.351
$ws.Range("L32").Value2 = 18907.615
$ws.Range("M32").Value2 = -5418.351
$ws.Range("N32").Value2 = -19481.615

# Row 61
$ws.Range("H61").Value2 = 2585142.2
$ws.Range("I61").Value2 = 2711149
$ws.Range("J61").Value2 = 2007
$ws.Range("K61").Value2 = 2711149
$ws.Range("L61").Value2 = 2007
$ws.Range("M61").Value2 = -2710937
$ws.Range("N61").Value2 = -2431

# Row 74
$ws.Range("H74").Value2 = 21434494
$ws.Range("I74").Value2 = 31579638
$ws.Range("J74").Value2 = 16967.555
$ws.Range("K74").Value2 = 31579638
$ws.Range("L74").Value2 = 16967.555
$ws.Range("M74").Value2 = -31578764
$ws.Range("N74").Value2 = -18715.555

# Row 77
$ws.Range("H77").Value2 = 21434494
$ws.Range("I77").Value2 = 31579638
$ws.Range("J77").Value2 = 16967.555
$ws.Range("K77").Value2 = 157898190
$ws.Range("L77").Value2 = 84837.77499999999
$ws.Range("M77").Value2 = -157893822
$ws.Range("N77").Value2 = -93573.77499999999

# Row 110
$ws.Range("H110").Value2 = 1315.625
$ws.Range("I110").Value2 = 1065.6666
$ws.Range("J110").Value2 = 2065.5
$ws.Range("K110").Value2 = 1065.6666
$ws.Range("L110").Value2 = 2065.5
$ws.Range("M110").Value2 = 979.3334
$ws.Range("N110").Value2 = -6155.5

# Row 122
$ws.Range("H122").Value2 = 2007
$ws.Range("I122").Value2 = 1610
$ws.Range("J122").Value2 = 2933.3333
$ws.Range("K122").Value2 = 4830
$ws.Range("L122").Value2 = 8799.999899999999
$ws.Range("M122").Value2 = -2380
$ws.Range("N122").Value2 = -13699.9999

# Row 132
$ws.Range("H132").Value2 = 808494.9399999999
$ws.Range("I132").Value2 = 886020.5
$ws.Range("J132").Value2 = 203795.6
$ws.Range("K132").Value2 = 2658061.5
$ws.Range("L132").Value2 = 611386.8
$ws.Range("M132").Value2 = -2655531.5
$ws.Range("N132").Value2 = -616446.8

# Row 136
$ws.Range("H136").Value2 = 2585142.2
$ws.Range("I136").Value2 = 2711149
$ws.Range("J136").Value2 = 2007
$ws.Range("K136").Value2 = 8133447
$ws.Range("L136").Value2 = 6021
$ws.Range("M136").Value2 = -8130897
$ws.Range("N136").Value2 = -11121

$ws = $wb.Worksheets.Item("BSM")
# Row 55
$ws.Range("H55").Value2 = 45291.668
$ws.Range("J55").Value2 = 45291.668
$ws.Range("L55").Value2 = 45291.668
$ws.Range("N55").Value2 = -45837.668

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 1889.6207
$ws.Range("I31").Value2 = 1403.579
$ws.Range("J31").Value2 = 2813.1
$ws.Range("K31").Value2 = 1403.579
$ws.Range("L31").Value2 = 2813.1
$ws.Range("M31").Value2 = -1108.579
$ws.Range("N31").Value2 = -3403.1

# Row 34
$ws.Range("H34").Value2 = 1889.6207
$ws.Range("I34").Value2 = 1403.579
$ws.Range("J34").Value2 = 2813.1
$ws.Range("K34").Value2 = 1403.579
$ws.Range("L34").Value2 = 2813.1
$ws.Range("M34").Value2 = -1201.579
$ws.Range("N34").Value2 = -3217.1

# Row 107
$ws.Range("H107").Value2 = 1894231.9
$ws.Range("I107").Value2 = 3472379.2
$ws.Range("J107").Value2 = 454.9
$ws.Range("K107").Value2 = 3472379.2
$ws.Range("L107").Value2 = 454.9
$ws.Range("M107").Value2 = -3470459.2
$ws.Range("N107").Value2 = -4294.9

# Row 109
$ws.Range("H109").Value2 = 48000
$ws.Range("J109").Value2 = 48000
$ws.Range("L109").Value2 = 48000
$ws.Range("N109").Value2 = -50080

# Row 134
$ws.Range("H134").Value2 = 2130.9688
$ws.Range("I134").Value2 = 2024.0344
$ws.Range("J134").Value2 = 3164.6667
$ws.Range("K134").Value2 = 6072.1032
$ws.Range("L134").Value2 = 9494.000100000001
$ws.Range("M134").Value2 = -3537.1032
$ws.Range("N134").Value2 = -14564.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value2 = 749.46
$ws.Range("I131").Value2 = 543.8333
$ws.Range("J131").Value2 = 777.5
$ws.Range("K131").Value2 = 1631.4999
$ws.Range("L131").Value2 = 2332.5
$ws.Range("M131").Value2 = 3408.5001
$ws.Range("N131").Value2 = -12412.5

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value2 = 2012.375
$ws.Range("I122").Value2 = 1638.8
$ws.Range("J122").Value2 = 2635
$ws.Range("K122").Value2 = 4916.4
$ws.Range("L122").Value2 = 7905
$ws.Range("M122").Value2 = -2466.4
$ws.Range("N122").Value2 = -12805

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value2 = 4977.3335
$ws.Range("I132").Value2 = 5191.75
$ws.Range("J132").Value2 = 3776.6
$ws.Range("K132").Value2 = 15575.25
$ws.Range("L132").Value2 = 11329.8
$ws.Range("M132").Value2 = -13045.25
$ws.Range("N132").Value2 = -16389.8

# Row 136
$ws.Range("H136").Value2 = 5879.375
$ws.Range("I136").Value2 = 8652
$ws.Range("J136").Value2 = 1258.3334
$ws.Range("K136").Value2 = 25956
$ws.Range("L136").Value2 = 3775.0002
$ws.Range("M136").Value2 = -23406
$ws.Range("N136").Value2 = -8875.0002
